# Update CLO grade-sheet scores: set the raw sub-objective scores (column B)
# on the "Grade Sheet" worksheet. The dependent ROUND(AVERAGE(...)) and
# count/letter-grade formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grade Sheet")

$ws.Range("B7").Value = 4
$ws.Range("B12").Value = 4
$ws.Range("B15").Value = 4
$ws.Range("B16").Value = 4
$ws.Range("B19").Value = 4
$ws.Range("B20").Value = 4
$ws.Range("B21").Value = 2
$ws.Range("B23").Value = 4
$ws.Range("B24").Value = 4
$ws.Range("B25").Value = 4
$ws.Range("B27").Value = 4
$ws.Range("B28").Value = 4
$ws.Range("B29").Value = 4
$ws.Range("B30").Value = 2
$ws.Range("B32").Value = 4
$ws.Range("B33").Value = 4
$ws.Range("B34").Value = 4
$ws.Range("B40").Value = 3
$ws.Range("B47").Value = 4
$ws.Range("B50").Value = 4
$ws.Range("B51").Value = 4

# Page setup: force portrait orientation (adds <pageSetup .../> on save).
$ws.PageSetup.Orientation = 1

# Move the window/selection to where the grader left off (row 21 area),
# with B40 as the active cell.
$ws.Activate()
$ws.Range("B40").Select()
